# Append two new data rows (19 and 20) for participant "fsgr7y5o",
# and widen the "nutrients" column (E) from 95 to 98 characters wide.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19
$ws.Range("A19").Value = "fsgr7y5o"
$ws.Range("B19").Value = "Training phase"
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = "['Purple', 'Orange', 'Orange', 'Orange', 'Green', 'Purple']"
$ws.Range("E19").Value = "[['Red', 'Red'], ['Blue', ''], ['Yellow', ''], ['Red', 'Blue'], ['Yellow', ''], ['Blue', 'Red']]"

# Row 20
$ws.Range("A20").Value = "fsgr7y5o"
$ws.Range("B20").Value = "Training phase"
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = "['Green', 'Green', 'Purple', 'Orange', 'Purple', 'Purple']"
$ws.Range("E20").Value = "[['Blue', ''], ['Red', ''], ['Red', ''], ['Yellow', ''], ['Yellow', ''], ['Blue', '']]"

# Widen column E (nutrients) from 95 to 98 raw OOXML width units.
# Excel's ColumnWidth property is offset from the stored XML width by
# the workbook's default "extra" padding (~0.8333 here), so add that
# back on to land exactly on 98 after save.
$ws.Columns.Item(5).ColumnWidth = 98 - 0.8333333333333334
